$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the section header rows (sex / mother's education / wealth quintile)
# Row 14 - "by sex" headers (Kyrgyz / Russian / English)
$ws.Range("A14").Value = "Жынысы боюнча"
$ws.Range("B14").Value = "По полу"
$ws.Range("C14").Value = "By sex"

# Row 17 - "education of mother" headers (Kyrgyz / Russian / English)
$ws.Range("A17").Value = "Энесинин билими "
$ws.Range("B17").Value = "Образование матери "
$ws.Range("C17").Value = "Education of mother"

# Row 23 - "wealth quintile" headers (Russian / English; Kyrgyz text unchanged)
$ws.Range("B23").Value = "Квинтиль по индексу благосостояния"
$ws.Range("C23").Value = "Wealth quintile"
